# "Generate Report for Handback"
# Updates the localization-status workbook to reflect a completed handback:
#   - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#     (this text is shared by Overview!E2/F2 and the zh-cn/de-de Status columns)
#   - zh-cn and de-de sheets get their "Latest Target File" (I2) and
#     "Latest Handback File" (J2) populated, with a hyperlink added on I2
#   - "Latest Handback DateTime" (K2) is stamped with the handback timestamp
#     (different timestamps for zh-cn vs de-de)
#   - A handful of columns are widened to better fit the newly-populated data

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$mdFileName = "285d2b4f-c17a-4342-84ac-1c63d0f35aa6.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/058f1c5ce4715f0482d3e340d18560d644410305/e2e/285d2b4f-c17a-4342-84ac-1c63d0f35aa6.md"
$zhCnXlf = "285d2b4f-c17a-4342-84ac-1c63d0f35aa6.625348eb0cedbb4f26c27554be30eed7f36c7f65.zh-cn.xlf"
$deDeXlf = "285d2b4f-c17a-4342-84ac-1c63d0f35aa6.625348eb0cedbb4f26c27554be30eed7f36c7f65.de-de.xlf"
$newStatus = "Handed back: in sync with en-US"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# Update every cell that currently shows the old status so the shared string
# text itself is replaced everywhere it is used.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- zh-cn: fill in Latest Target File / Latest Handback File / Handback DateTime ---
$wsZhCn.Range("I2").Value = $mdFileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null
$f = $wsZhCn.Range("I2").Font
$f.Underline = $true
$f.Color = 15570276
$wsZhCn.Range("J2").Value = $zhCnXlf
$wsZhCn.Range("K2").Value = "2016-09-04 17:04:46"

# --- de-de: fill in Latest Target File / Latest Handback File / Handback DateTime ---
$wsDeDe.Range("I2").Value = $mdFileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null
$f2 = $wsDeDe.Range("I2").Font
$f2.Underline = $true
$f2.Color = 15570276
$wsDeDe.Range("J2").Value = $deDeXlf
$wsDeDe.Range("K2").Value = "2016-09-04 17:04:53"

# --- Column widening to fit the now-populated columns ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

$wsZhCn.Columns.Item(3).ColumnWidth = 29.14
$wsZhCn.Columns.Item(9).ColumnWidth = 39.16
$wsZhCn.Columns.Item(10).ColumnWidth = 39.16

$wsDeDe.Columns.Item(3).ColumnWidth = 29.14
$wsDeDe.Columns.Item(9).ColumnWidth = 39.16
$wsDeDe.Columns.Item(10).ColumnWidth = 39.16
